$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert new row at position 3 (shifts old rows 3-6 down to 4-7, including column A)
$ws2.Rows("3:3").Insert()

# Clear the unwanted C3/H3 placeholder cells entirely
$ws2.Range("C3").Clear()
$ws2.Range("H3").Clear()

# Fix formatting of the new row 3 cells to match surrounding style
$ws2.Range("A4").Copy()
$ws2.Range("A3").PasteSpecial(-4122)
$ws2.Range("E5").Copy()
$ws2.Range("E3").PasteSpecial(-4122)
$ws2.Range("G4").Copy()
$ws2.Range("G3").PasteSpecial(-4122)
$ws2.Range("H6").Copy()
$ws2.Range("D3").PasteSpecial(-4122)

# Restore A column step numbers (Insert shifted col A down too; put back the sequential numbering)
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5
$ws2.Range("A7").Value = 6

# Set values for new row 3 (the inserted action)
# Order matters for shared-string table append order: gmail-xpath, compto, Hola
$ws2.Range("F3").Value = "//a[contains(text(),'Gmail')]"
$ws2.Range("B3").Value = "compto"
$ws2.Range("D3").Value = "Hola"
$ws2.Range("E3").Value = "xpath"
$ws2.Range("G3").Value = $true

# Update selections / active sheet: TC2 becomes the active (visible) tab,
# TC1 selection moves to D1, TC2 selection moves to G3.
$ws1.Activate()
$ws1.Range("D1").Select() | Out-Null

$ws2.Activate()
$ws2.Range("G3").Select() | Out-Null

Write-Host "done"
